$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 'RT @gDesFaits: #BIDEN sait bien s''entourer... #pedophilie
Jerry Harris, star de l''émission Netflix « Cheer » et ancien substitut de la cam…'
$ws.Range("C2").Value = 1339229134069915904
$ws.Range("D2").Value = 140
$ws.Range("E2").Value = 44181.63982638889
$ws.Range("F2").Value = 'fr'
$ws.Range("G2").Value = 1077271244964315008
$ws.Range("H2").Value = 'botduSEXE'
$ws.Range("I2").Value = 1850
$ws.Range("J2").Value = ''
$ws.Range("K2").Value = 'RT @gDesFaits: #BIDEN knows how to surround himself well ... #pedophiliaJerry Harris, star of the Netflix show "Cheer" and former cam substitute ... '
$ws.Range("L2").Value = ' #BIDEN sait bien s''entourer... #pedophilie  Jerry Harris, star de l''émission Netflix « Cheer » et ancien substitut de la cam…'
$ws.Range("M2").Value = ' #BIDEN knows how to surround himself well ... #pedophiliaJerry Harris, star of the Netflix show "Cheer" and former cam substitute ... '
$ws.Range("N2").Value = 'RT @gDesFaits: #BIDEN knows how to surround himself well ... #pedophiliaJerry Harris, star of the Netflix show "Cheer" and former cam substitute ... '

# Row 3
$ws.Range("B3").Value = 'Abusing #Republican Leaders #Media #Trump #Facebook America Law Government Sports Entertainment health Education Am… https://t.co/2uhi1nNOjP'
$ws.Range("C3").Value = 1339229096069349120
$ws.Range("D3").Value = 140
$ws.Range("E3").Value = 44181.63972222222
$ws.Range("F3").Value = 'en'
$ws.Range("G3").Value = 1352572483
$ws.Range("H3").Value = 'sethrow991'
$ws.Range("I3").Value = 162
$ws.Range("J3").Value = ''
$ws.Range("K3").Value = 'Abusing #Republican Leaders #Media #Trump #Facebook America Law Government Sports Entertainment health Education Am… https://t.co/2uhi1nNOjP '
$ws.Range("L3").Value = 'Abusing #Republican Leaders #Media #Trump #Facebook America Law Government Sports Entertainment health Education Am… '
$ws.Range("M3").Value = 'Abusing #Republican Leaders #Media #Trump #Facebook America Law Government Sports Entertainment health Education Am…  '
$ws.Range("N3").Value = 'Abusing #Republican Leaders #Media #Trump #Facebook America Law Government Sports Entertainment health Education Am…  '

# Row 4
$ws.Range("B4").Value = 'RT @in_pubs: If you need a pick me up today....
Just think....
At this very moment @realDonaldTrump and his lawyers. 
Are scrambling to…'
$ws.Range("C4").Value = 1339229085529231104
$ws.Range("D4").Value = 139
$ws.Range("E4").Value = 44181.63969907408
$ws.Range("F4").Value = 'en'
$ws.Range("G4").Value = 822589145335955456
$ws.Range("H4").Value = 'Loiskane1202'
$ws.Range("I4").Value = 154
$ws.Range("J4").Value = ''
$ws.Range("K4").Value = 'RT @in_pubs: If you need a pick me up today....Just think....At this very moment @realDonaldTrump and his lawyers. Are scrambling to… '
$ws.Range("L4").Value = ' If you need a pick me up today....  Just think....  At this very moment  and his lawyers.   Are scrambling to…'
$ws.Range("M4").Value = ' If you need a pick me up today....Just think....At this very moment  and his lawyers. Are scrambling to… '
$ws.Range("N4").Value = 'RT @in_pubs: If you need a pick me up today....Just think....At this very moment @realDonaldTrump and his lawyers. Are scrambling to… '

# Row 5
$ws.Range("B5").Value = 'Oups another #freakout! #trump #fail realDonaldTrump: Perhaps the biggest difference between 2016 and 2020 is… https://t.co/LjFKYUVHx9'
$ws.Range("C5").Value = 1339229062091443968
$ws.Range("D5").Value = 134
$ws.Range("E5").Value = 44181.63962962963
$ws.Range("F5").Value = 'en'
$ws.Range("G5").Value = 839482284759724032
$ws.Range("H5").Value = 'trumpfreakout'
$ws.Range("I5").Value = 67
$ws.Range("J5").Value = ''
$ws.Range("K5").Value = 'Oups another #freakout! #trump #fail realDonaldTrump: Perhaps the biggest difference between 2016 and 2020 is… https://t.co/LjFKYUVHx9 '
$ws.Range("L5").Value = 'Oups another #freakout! #trump #fail realDonaldTrump: Perhaps the biggest difference between 2016 and 2020 is… '
$ws.Range("M5").Value = 'Oups another #freakout! #trump #fail realDonaldTrump: Perhaps the biggest difference between 2016 and 2020 is…  '
$ws.Range("N5").Value = 'Oups another #freakout! #trump #fail realDonaldTrump: Perhaps the biggest difference between 2016 and 2020 is…  '

# Row 6
$ws.Range("B6").Value = '@ThomTillis @Perduesenate @KLoeffler Like the other #GOPTraitors, #Loeffler & #Perdue stood silently by as #Trump d… https://t.co/tPSF80ZcfZ'
$ws.Range("C6").Value = 1339229040096456960
$ws.Range("D6").Value = 144
$ws.Range("E6").Value = 44181.63957175926
$ws.Range("F6").Value = 'en'
$ws.Range("G6").Value = 416363599
$ws.Range("H6").Value = 'GeneKelsey'
$ws.Range("I6").Value = 113
$ws.Range("J6").Value = 'Florida'
$ws.Range("K6").Value = '@ThomTillis @Perduesenate @KLoeffler Like the other #GOPTraitors, #Loeffler & #Perdue stood silently by as #Trump d… https://t.co/tPSF80ZcfZ '
$ws.Range("L6").Value = '   Like the other #GOPTraitors, #Loeffler & #Perdue stood silently by as #Trump d… '
$ws.Range("M6").Value = '   Like the other #GOPTraitors, #Loeffler & #Perdue stood silently by as #Trump d…  '
$ws.Range("N6").Value = '@ThomTillis @Perduesenate @KLoeffler Like the other #GOPTraitors, #Loeffler & #Perdue stood silently by as #Trump d…  '

# Row 7
$ws.Range("B7").Value = 'IN 2020, THE #TRUMP ADMINISTRATION DECLARED WAR ON DANCING TEENS
https://t.co/dumJMhzp02'
$ws.Range("C7").Value = 1339229035797351936
$ws.Range("D7").Value = 88
$ws.Range("E7").Value = 44181.63956018518
$ws.Range("F7").Value = 'en'
$ws.Range("G7").Value = 90272103
$ws.Range("H7").Value = 'gezgintrk'
$ws.Range("I7").Value = 12308
$ws.Range("J7").Value = 'Turkey / İstanbul'
$ws.Range("K7").Value = 'IN 2020, THE #TRUMP ADMINISTRATION DECLARED WAR ON DANCING TEENShttps://t.co/dumJMhzp02 '
$ws.Range("L7").Value = 'IN 2020, THE #TRUMP ADMINISTRATION DECLARED WAR ON DANCING TEENS '
$ws.Range("M7").Value = 'IN 2020, THE #TRUMP ADMINISTRATION DECLARED WAR ON DANCING TEENS '
$ws.Range("N7").Value = 'IN 2020, THE #TRUMP ADMINISTRATION DECLARED WAR ON DANCING TEENS '

# Row 8
$ws.Range("B8").Value = 'RT @gDesFaits: #BIDEN sait bien s''entourer... #pedophilie
Jerry Harris, star de l''émission Netflix « Cheer » et ancien substitut de la cam…'
$ws.Range("C8").Value = 1339229026062393088
$ws.Range("D8").Value = 140
$ws.Range("E8").Value = 44181.63953703704
$ws.Range("F8").Value = 'fr'
$ws.Range("G8").Value = 716571320574676992
$ws.Range("H8").Value = 'JePPauwels'
$ws.Range("I8").Value = 154
$ws.Range("J8").Value = 'Hainaut, Belgique'
$ws.Range("K8").Value = 'RT @gDesFaits: #BIDEN knows how to surround himself well ... #pedophiliaJerry Harris, star of the Netflix show "Cheer" and former cam substitute ... '
$ws.Range("L8").Value = ' #BIDEN sait bien s''entourer... #pedophilie  Jerry Harris, star de l''émission Netflix « Cheer » et ancien substitut de la cam…'
$ws.Range("M8").Value = ' #BIDEN knows how to surround himself well ... #pedophiliaJerry Harris, star of the Netflix show "Cheer" and former cam substitute ... '
$ws.Range("N8").Value = 'RT @gDesFaits: #BIDEN knows how to surround himself well ... #pedophiliaJerry Harris, star of the Netflix show "Cheer" and former cam substitute ... '

# Row 9
$ws.Range("B9").Value = 'RT @MarteauOlivier: Quand #Obama a été réélu en 2012, il avait perdu 3 millions de voix par rapport à son élection en 2008.
#Trump lui a ga…'
$ws.Range("C9").Value = 1339228998992358912
$ws.Range("D9").Value = 140
$ws.Range("E9").Value = 44181.63945601852
$ws.Range("F9").Value = 'fr'
$ws.Range("G9").Value = 3138038768
$ws.Range("H9").Value = 'Randy64_fr'
$ws.Range("I9").Value = 331
$ws.Range("J9").Value = ''
$ws.Range("K9").Value = 'RT @MarteauOlivier: When #Obama was re-elected in 2012, he had lost 3 million votes compared to his election in 2008. # Trump gave him ... '
$ws.Range("L9").Value = ' Quand #Obama a été réélu en 2012, il avait perdu 3 millions de voix par rapport à son élection en 2008. #Trump lui a ga…'
$ws.Range("M9").Value = ' When #Obama was re-elected in 2012, he had lost 3 million votes compared to his election in 2008. # Trump gave him ... '
$ws.Range("N9").Value = 'RT @MarteauOlivier: When #Obama was re-elected in 2012, he had lost 3 million votes compared to his election in 2008. # Trump gave him ... '

# Row 10
$ws.Range("B10").Value = 'RT @Susan10515068: Trump fucking lost and there is nothing you can do...
#TrumpTheFool 
#Trump 
#PsychoTrump https://t.co/wWvLaJzEig'
$ws.Range("C10").Value = 1339228997310406912
$ws.Range("D10").Value = 132
$ws.Range("E10").Value = 44181.63945601852
$ws.Range("F10").Value = 'en'
$ws.Range("G10").Value = 900977000
$ws.Range("H10").Value = 'Cradd4Teresa'
$ws.Range("I10").Value = 29
$ws.Range("J10").Value = ''
$ws.Range("K10").Value = 'RT @Susan10515068: Trump fucking lost and there is nothing you can do...#TrumpTheFool #Trump #PsychoTrump https://t.co/wWvLaJzEig '
$ws.Range("L10").Value = ' Trump fucking lost and there is nothing you can do... #TrumpTheFool  #Trump  #PsychoTrump '
$ws.Range("M10").Value = ' Trump fucking lost and there is nothing you can do...#TrumpTheFool #Trump #PsychoTrump  '
$ws.Range("N10").Value = 'RT @Susan10515068: Trump fucking lost and there is nothing you can do...#TrumpTheFool #Trump #PsychoTrump  '

# Row 11
$ws.Range("B11").Value = '@realDonaldTrump @FoxNews #Trump is irrelevant. 
Ignore the insane ramblings of #PsychoTrump #25th 
#Georgia deser… https://t.co/RYfE4nwE8B'
$ws.Range("C11").Value = 1339228980302536960
$ws.Range("D11").Value = 140
$ws.Range("E11").Value = 44181.63940972222
$ws.Range("F11").Value = 'en'
$ws.Range("G11").Value = 1173506442491518976
$ws.Range("H11").Value = 'LeeSaunders72'
$ws.Range("I11").Value = 1179
$ws.Range("J11").Value = 'Leeds/Manchester/London, UK'
$ws.Range("K11").Value = '@realDonaldTrump @FoxNews #Trump is irrelevant. Ignore the insane ramblings of #PsychoTrump #25th #Georgia deser… https://t.co/RYfE4nwE8B '
$ws.Range("L11").Value = '  #Trump is irrelevant.  Ignore the insane ramblings of #PsychoTrump #25th   #Georgia deser… '
$ws.Range("M11").Value = '  #Trump is irrelevant. Ignore the insane ramblings of #PsychoTrump #25th #Georgia deser…  '
$ws.Range("N11").Value = '@realDonaldTrump @FoxNews #Trump is irrelevant. Ignore the insane ramblings of #PsychoTrump #25th #Georgia deser…  '
